# Applies the "AVO_001_0300_Exportar_Facturas" config update to Config.xlsx:
#  - Rename the process/queue value and the business-process-name value from
#    the old framework placeholders to the new process name.
#  - Add three new config rows (PDFOriginPath, PDFDestinationPath,
#    XMLDestinationPath) with their descriptions on the Settings sheet.
#  - Leave the cursor/selection the way the author apparently left it:
#    Constants!C17, Assets!A2, and finally Settings!A10 (so Settings ends
#    up the active tab/sheet).

$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# --- Settings sheet: rename the two placeholder values -------------------
$wsSettings.Range("B2").Value = "AVO_001_0300_Exportar_Facturas"
$wsSettings.Range("B5").Value = "AVO_001_0300_Exportar_Facturas"

# --- Settings sheet: new PDF/XML path configuration rows -----------------
$wsSettings.Range("A7").Value = "PDFOriginPath"
$wsSettings.Range("C7").Value = "Ruta donde se almacena el documento original de la factura"
$wsSettings.Range("C7").WrapText = $true

$wsSettings.Range("A8").Value = "PDFDestinationPath"
$wsSettings.Range("C8").Value = "Ruta donde se almacenará el documento PDF en los sistemas del cliente"

$wsSettings.Range("A9").Value = "XMLDestinationPath"
$wsSettings.Range("C9").Value = "Ruta donde se almacenará el documento XML generado"
$wsSettings.Range("C9").WrapText = $true

# --- Restore the per-sheet cursor positions, then land on Settings!A10 ---
$wsConstants.Activate() | Out-Null
$wsConstants.Range("C17").Select() | Out-Null

$wsAssets.Activate() | Out-Null
$wsAssets.Range("A2").Select() | Out-Null

$wsSettings.Activate() | Out-Null
$wsSettings.Range("A10").Select() | Out-Null
